# Apply NBO data updates to Sheet1.
# The edit refreshes the "Inactive Val51" orbital-interaction block (rows 23-34,
# columns F:I) with a newer round of NBO second-order perturbation values, and
# moves the saved selection/scroll position in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated orbital interaction data (columns F:I) for rows 23-34 ---
# Row layout: Orbital label, Interaction with LP1, LP2, LP3 (column J is the
# existing shared SUM formula and recalculates automatically).

$updates = @(
    @{ Row = 23; F = "BD*( 2) C232- O233"; G = 0.1;  H = 0.07; I = "n/a" },
    @{ Row = 24; F = "BD*( 1) C 44- H 64"; G = 0;    H = 0;    I = "n/a" },
    @{ Row = 25; F = "BD*( 1) N 41- H 65"; G = 9.16; H = 1.66; I = "n/a" },
    @{ Row = 26; F = "BD*( 1) N 69- H 84"; G = 4.49; H = 0.48; I = "n/a" },
    @{ Row = 27; F = "BD*( 2) C232- O233"; G = 0.23; H = 0.11; I = "n/a" },
    @{ Row = 28; F = "BD*( 1) C212- H140"; G = 0;    H = 0;    I = 0 },
    @{ Row = 29; F = "BD*( 1) N117- H141"; G = 7.31; H = 4.97; I = 9.82 },
    @{ Row = 30; F = "BD*( 1) N 5- H 30";  G = 7.2;  H = 2.34; I = 9.63 },
    @{ Row = 31; F = "BD*( 2) C232- O233"; G = 0.15; H = 0.09; I = 0 },
    @{ Row = 32; F = "BD*( 1) C205- H225"; G = 0;    H = 0;    I = 0 },
    @{ Row = 33; F = "BD*( 1) N202- H226"; G = 7.71; H = 3.36; I = 11.09 },
    @{ Row = 34; F = "BD*( 1) N 44- H 68"; G = 6.76; H = 0;    I = 8.37 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("F$r").Value = $u.F
    $ws.Range("G$r").Value = $u.G
    $ws.Range("H$r").Value = $u.H
    $ws.Range("I$r").Value = $u.I
}

# --- Update the sheet's saved selection/scroll position to reflect where the
#     author left off while reviewing the table (row 15 scrolled to top, with
#     N30 the active cell of the selection) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N30").Select()
